# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The sheet tracks a "last changed" date stamp in column C that gets
# bumped forward by one day (45179 -> 45180, serial date numbers) for
# every data row (rows 2-302).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 302; $row++) {
    $ws.Cells.Item($row, 3).Value = 45180
}
